$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '90.893.67'
Set-TextValue 2 5 '  +1.16%  '

# Row 3
Set-TextValue 3 4 '3.175.85'
Set-TextValue 3 5 '  +2.35%  '

# Row 4
Set-TextValue 4 5 '  +0.48%  '

# Row 5
Set-TextValue 5 4 '219.72'
Set-TextValue 5 5 '  +2.58%  '

# Row 6
Set-TextValue 6 4 '625.47'
Set-TextValue 6 5 '  +0.96%  '

# Row 7
Set-TextValue 7 5 '  +21.26%  '

# Row 8
Set-TextValue 8 4 '0.374'
Set-TextValue 8 5 '  +0.50%  '

# Row 9
Set-TextValue 9 5 '  +0.01%  '

# Row 10
Set-TextValue 10 4 '3.173.20'
Set-TextValue 10 5 '  +2.37%  '

# Row 11
Set-TextValue 11 4 '0.749'
Set-TextValue 11 5 '  +17.68%  '

# Row 12
Set-TextValue 12 4 '0.199'
Set-TextValue 12 5 '  +5.06%  '

# Row 13
Set-TextValue 13 4 '0.0000251'
Set-TextValue 13 5 '  +3.79%  '

# Row 14
Set-TextValue 14 4 '35.10'
Set-TextValue 14 5 '  +8.48%  '

# Row 15
Set-TextValue 15 4 '5.56'
Set-TextValue 15 5 '  +4.48%  '

# Row 16
Set-TextValue 16 4 '90.832.08'
Set-TextValue 16 5 '  +1.31%  '

# Row 17
Set-TextValue 17 4 '3.781.31'
Set-TextValue 17 5 '  +2.47%  '

# Row 18
Set-TextValue 18 4 '3.167.66'
Set-TextValue 18 5 '  +1.69%  '

# Row 19
Set-TextValue 19 4 '3.77'
Set-TextValue 19 5 '  +9.91%  '

# Row 20
Set-TextValue 20 4 '0.0000220'
Set-TextValue 20 5 '  +0.68%  '

# Row 21
Set-TextValue 21 4 '14.37'
Set-TextValue 21 5 '  +6.76%  '

# Row 22
Set-TextValue 22 4 '443.81'
Set-TextValue 22 5 '  +3.81%  '

# Row 23
Set-TextValue 23 5 '  +10.07%  '

# Row 24
Set-TextValue 24 5 '  +5.21%  '

# Row 25
Set-TextValue 25 5 '  +11.89%  '

# Row 26
Set-TextValue 26 4 '87.31'
Set-TextValue 26 5 '  +4.04%  '

# Row 27
Set-TextValue 27 4 '12.33'
Set-TextValue 27 5 '  +2.34%  '

# Row 28
Set-TextValue 28 4 '3.345.77'
Set-TextValue 28 5 '  +1.40%  '

# Row 29
Set-TextValue 29 4 '0.997'
Set-TextValue 29 5 '  -0.27%  '

# Row 30
Set-TextValue 30 4 '0.165'
Set-TextValue 30 5 '  +1.85%  '

# Row 31
Set-TextValue 31 4 '9.30'
Set-TextValue 31 5 '  +14.16%  '

# Row 32
Set-TextValue 32 4 '0.996'
Set-TextValue 32 5 '  -7.82%  '

# Row 33
Set-TextValue 33 4 '529.60'
Set-TextValue 33 5 '  +3.53%  '

# Row 34
Set-TextValue 34 4 '25.27'
Set-TextValue 34 5 '  +11.92%  '

# Row 35
Set-TextValue 35 4 '3.77'
Set-TextValue 35 5 '  +1.95%  '

# Row 36
Set-TextValue 36 2 'Kaspa'
Set-TextValue 36 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 36 4 '0.146'
Set-TextValue 36 5 '  +12.26%  '

# Row 37
Set-TextValue 37 2 'RenderToken'
Set-TextValue 37 3 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 37 4 '7.05'
Set-TextValue 37 5 '  +4.94%  '

# Row 38
Set-TextValue 38 4 '1.91'
Set-TextValue 38 5 '  +6.20%  '

# Row 39
Set-TextValue 39 4 '1.31'
Set-TextValue 39 5 '  +4.77%  '

# Row 40
Set-TextValue 40 5 '  -0.28%  '

# Row 41
Set-TextValue 41 5 '  +13.08%  '

# Row 42
Set-TextValue 42 2 'Hedera'
Set-TextValue 42 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 42 4 '0.0857'
Set-TextValue 42 5 '  +21.84%  '

# Row 43
Set-TextValue 43 2 'FirstDigitalUSD'
Set-TextValue 43 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 43 4 '1.00'
Set-TextValue 43 5 '  -0.13%  '

# Row 44
Set-TextValue 44 4 '0.411'
Set-TextValue 44 5 '  +11.62%  '

# Row 45
Set-TextValue 45 4 '1.96'
Set-TextValue 45 5 '  +6.09%  '

# Row 46
Set-TextValue 46 5 '  +0.02%  '

# Row 47
Set-TextValue 47 4 '148.32'
Set-TextValue 47 5 '  +1.63%  '

# Row 48
Set-TextValue 48 5 '  +9.88%  '

# Row 49
Set-TextValue 49 4 '44.47'
Set-TextValue 49 5 '  +1.92%  '

# Row 50
Set-TextValue 50 2 'Filecoin'
Set-TextValue 50 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 50 4 '4.40'
Set-TextValue 50 5 '  +10.13%  '

# Row 51
Set-TextValue 51 2 'Aave'
Set-TextValue 51 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 51 4 '171.62'
Set-TextValue 51 5 '  +7.68%  '
